# Updates cryptos list with refreshed price / 1h-volume figures, plus three
# row-content swaps (SuiNetwork<->PEPE at rows 19-20, and the 3-way rotation
# Fetch.AI -> Kaspa -> PancakeSwap -> Fetch.AI across rows 36-38).
#
# Numeric-looking price strings (single decimal point, e.g. "0.997", "4.05")
# are written with a temporary Text number format so Excel keeps the exact
# literal text (incl. trailing zeros) instead of auto-coercing to a number;
# ClearFormats() immediately restores the cell's original (default) style.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "88.993.79"
$ws.Range("E2").Value = "  -3.72%  "
$ws.Range("D3").Value = "3.161.59"
$ws.Range("E3").Value = "  -3.83%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.997"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  -0.44%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "216.46"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.07%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "635.66"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +0.17%  "
$ws.Range("E7").Value = "  -3.00%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.729"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +1.20%  "
$ws.Range("E9").Value = "  +0.00%  "
$ws.Range("D10").Value = "3.158.06"
$ws.Range("E10").Value = "  -3.72%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.558"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -5.11%  "
$ws.Range("E12").Value = "  -0.38%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000253"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -3.78%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.32"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -0.95%  "
$ws.Range("D15").Value = "88.494.44"
$ws.Range("E15").Value = "  -3.88%  "
$ws.Range("D16").Value = "3.723.39"
$ws.Range("E16").Value = "  -4.33%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "32.65"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -4.47%  "
$ws.Range("D18").Value = "3.144.45"
$ws.Range("E18").Value = "  -4.02%  "
$ws.Range("B19").Value = "PEPE"
$ws.Range("C19").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.0000231"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +22.50%  "
$ws.Range("B20").Value = "SuiNetwork"
$ws.Range("C20").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "3.40"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +2.14%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.33"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -5.32%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "428.28"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -2.45%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.44"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -5.38%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "4.92"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -7.02%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "5.41"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +0.44%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.62"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -5.44%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "80.52"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +4.99%  "
$ws.Range("D28").Value = "3.301.39"
$ws.Range("E28").Value = "  -6.05%  "
$ws.Range("E29").Value = "  +0.02%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.161"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -11.24%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.995"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -0.23%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.05"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +10.20%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "8.27"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -6.05%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "515.87"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -7.92%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "7.16"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +0.16%  "
$ws.Range("B36").Value = "Kaspa"
$ws.Range("C36").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.144"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +9.19%  "
$ws.Range("B37").Value = "Fetch.AI"
$ws.Range("C37").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.32"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +2.81%  "
$ws.Range("B38").Value = "PancakeSwap"
$ws.Range("C38").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.84"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -4.70%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "22.04"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -2.99%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "22.24"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -1.02%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.998"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -0.14%  "
$ws.Range("E42").Value = "  -0.13%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.88"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -6.35%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.367"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -7.21%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "145.73"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -4.20%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "43.77"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -0.44%  "
$ws.Range("E47").Value = "  -2.96%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "167.19"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -7.83%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.730"
$ws.Range("D49").ClearFormats()
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "24.83"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -0.68%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.20"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -7.05%  "
